$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 74 (B74): status changes from "no comenzado" to "terminado"
$ws.Range("B74").Value = "terminado"

# New row 83: new task "costos de articulos calcular aquellos q estan formulados"
# with status "no comenzado"
$ws.Range("A83").Value = "costos de articulos calcular aquellos q estan formulados"
$ws.Range("B83").Value = "no comenzado"

# Update the active selection to match the author's final cursor position
$ws.Range("B74").Select()
